# Update betting odds on Sheet1 to reflect the latest FlashScore snapshot.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 6
$ws.Range("G6").Value = 2.3
$ws.Range("I6").Value = 3.7
$ws.Range("J6").Value = 3.25
$ws.Range("L6").Value = 4.5
$ws.Range("M6").Value = 1.17
$ws.Range("N6").Value = 5
$ws.Range("AH6").Value = 7
$ws.Range("AY6").Value = 23
$ws.Range("BA6").Value = 81

# Row 7
$ws.Range("Q7").Value = 2.15
$ws.Range("R7").Value = 1.67
$ws.Range("AW7").Value = 126

# Row 8
$ws.Range("K8").Value = 1.95
$ws.Range("Q8").Value = 2.5
$ws.Range("R8").Value = 1.5

# Row 9
$ws.Range("K9").Value = 1.8

# Row 17
$ws.Range("S17").Value = 1.37

# Row 47
$ws.Range("G47").Value = 1.26
$ws.Range("H47").Value = 5.2
$ws.Range("I47").Value = 10
$ws.Range("K47").Value = 2.47
$ws.Range("L47").Value = 8
$ws.Range("Q47").Value = 1.53
$ws.Range("R47").Value = 2.18
$ws.Range("S47").Value = 1.28
$ws.Range("T47").Value = 3.34
$ws.Range("U47").Value = 2
$ws.Range("Z47").Value = 7.5
$ws.Range("AD47").Value = 10.75
$ws.Range("AH47").Value = 27
$ws.Range("AI47").Value = 80
$ws.Range("AJ47").Value = 32
$ws.Range("AK47").Value = 300
$ws.Range("AO47").Value = 5.4
$ws.Range("AQ47").Value = 13.5
$ws.Range("AT47").Value = 3.1
$ws.Range("AX47").Value = 10
$ws.Range("AY47").Value = 60
$ws.Range("AZ47").Value = 55
$ws.Range("BA47").Value = 450
$ws.Range("BB47").Value = 450

# Row 62
$ws.Range("G62").Value = 2.5
$ws.Range("I62").Value = 2.7
$ws.Range("L62").Value = 3.6
$ws.Range("U62").Value = 2
$ws.Range("V62").Value = 1.73
$ws.Range("W62").Value = 7
$ws.Range("Y62").Value = 11
$ws.Range("Z62").Value = 26
$ws.Range("AJ62").Value = 11
$ws.Range("AL62").Value = 26
$ws.Range("AN62").Value = 4.5
$ws.Range("AY62").Value = 17
$ws.Range("AZ62").Value = 29
$ws.Range("BA62").Value = 51
$ws.Range("BC62").Value = 251

# Row 63
$ws.Range("G63").Value = 2.3
$ws.Range("J63").Value = 3.1
$ws.Range("M63").Value = 1.06
$ws.Range("N63").Value = 9.5
$ws.Range("O63").Value = 1.33
$ws.Range("P63").Value = 3.25
$ws.Range("Q63").Value = 2.08
$ws.Range("R63").Value = 1.73
$ws.Range("Y63").Value = 9.5
$ws.Range("AB63").Value = 29
$ws.Range("AC63").Value = 9.5
$ws.Range("AH63").Value = 8.5
$ws.Range("AN63").Value = 4.33
$ws.Range("AO63").Value = 13
$ws.Range("AP63").Value = 23
$ws.Range("AQ63").Value = 41
